$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Move existing comments out of the way (right to left) before inserting columns,
#        since Excel COM does not auto-shift cell comments when columns are inserted. ---
function Move-CommentText {
    param($fromAddr, $toAddr)
    $fromCell = $ws.Range($fromAddr)
    $comment = $fromCell.Comment
    if ($comment) {
        $txt = $comment.Text()
        $comment.Delete()
        $ws.Range($toAddr).AddComment($txt)
    }
}

# --- 2. Insert the two new columns J:K (shifts J->L, K->M, L->N, M->O, N->P data) ---
$ws.Columns("J:K").Insert()

# --- 3. Re-home the old comments (which stayed anchored to their original cell refs)
#        to their new shifted locations. Go right-to-left to avoid clobbering. ---
Move-CommentText "N1" "P1"
Move-CommentText "M1" "O1"
Move-CommentText "L1" "N1"
Move-CommentText "K1" "M1"
Move-CommentText "J1" "L1"

# --- 4. Add the two new comments on the newly inserted J1/K1 header cells ---
$ws.Range("J1").AddComment("No aparece en excel dicotom; el criterio debería ser 1 si la respuesta es menor a 40")
$ws.Range("K1").AddComment("Criterio de preg_edad, donde es 1 si la respuesta es <= 40")

# --- 5. Set the new header names ---
$ws.Range("J1").Value = "preg_edad"
$ws.Range("K1").Value = "crit_edad"

# --- 6. Set the new data values ---
$ws.Range("J2").Value = 30
$ws.Range("J3").Value = 60

$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "1"

$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "0"
